$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.831836938858032
$ws.Range("B1").Value = 1.966942667961121
$ws.Range("C1").Value = 1.949584364891052
$ws.Range("D1").Value = 2.331913948059082
$ws.Range("E1").Value = 3.081932067871094
